$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-12-10 Tuesday"; new = "2024-12-11 Wednesday"},
    @{old = "94×29=2726"; new = "93×96=8928"},
    @{old = "83×96=7968"; new = "72×40=2880"},
    @{old = "96×88=8448"; new = "92×76=6992"},
    @{old = "28×33=924"; new = "38×33=1254"},
    @{old = "81×85=6885"; new = "13×84=1092"},
    @{old = "57×30=1710"; new = "20×56=1120"},
    @{old = "36×51=1836"; new = "24×99=2376"},
    @{old = "71×51=3621"; new = "84×83=6972"},
    @{old = "17×44=748"; new = "27×30=810"},
    @{old = "92×87=8004"; new = "29×88=2552"},
    @{old = "61×86=5246"; new = "47×77=3619"},
    @{old = "24×23=552"; new = "26×65=1690"},
    @{old = "43×72=3096"; new = "62×62=3844"},
    @{old = "19×23=437"; new = "79×55=4345"},
    @{old = "40×44=1760"; new = "39×78=3042"},
    @{old = "88×59=5192"; new = "24×99=2376"},
    @{old = "81×49=3969"; new = "64×61=3904"},
    @{old = "39×65=2535"; new = "95×14=1330"},
    @{old = "12×15=180"; new = "21×55=1155"},
    @{old = "36×69=2484"; new = "14×88=1232"},
    @{old = "66×46=3036"; new = "36×83=2988"},
    @{old = "44×18=792"; new = "73×68=4964"},
    @{old = "47×38=1786"; new = "60×72=4320"},
    @{old = "52×11=572"; new = "93×76=7068"},
    @{old = "63×67=4221"; new = "47×73=3431"}
)

foreach ($rep in $replacements) {
    $d.Content.Find.Execute($rep.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $rep.new, 2)
}
